$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell F1 "time_taken", copying the style/format of E1 (the
# neighboring header cell) so it keeps the same bold/bordered header look.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Fill in the time_taken values for each data row (F2:F15)
$timestamps = @(
    "2021-10-05 13:38:28.955999",
    "2021-10-05 13:38:28.956011",
    "2021-10-05 13:38:28.956015",
    "2021-10-05 13:38:28.956017",
    "2021-10-05 13:38:28.956020",
    "2021-10-05 13:38:28.956023",
    "2021-10-05 13:38:28.956025",
    "2021-10-05 13:38:28.956028",
    "2021-10-05 13:38:28.956030",
    "2021-10-05 13:38:28.956033",
    "2021-10-05 13:38:28.956035",
    "2021-10-05 13:38:28.956038",
    "2021-10-05 13:38:28.956040",
    "2021-10-05 13:38:28.956043"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
